$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 313.81818
$ws.Range("J17").Value = 313.81818
$ws.Range("L17").Value = 941.45454
$ws.Range("N17").Value = -1277.45454

$ws.Range("H28").Value = 571.4
$ws.Range("I28").Value = 548.26086
$ws.Range("K28").Value = 548.26086
$ws.Range("M28").Value = -63.26085999999998

$ws.Range("H86").Value = 2804.8572
$ws.Range("I86").Value = 1350.875
$ws.Range("J86").Value = 4743.5
$ws.Range("K86").Value = 1350.875
$ws.Range("L86").Value = 4743.5
$ws.Range("M86").Value = -227.875
$ws.Range("N86").Value = -6989.5

$ws.Range("H89").Value = 2804.8572
$ws.Range("I89").Value = 1350.875
$ws.Range("J89").Value = 4743.5
$ws.Range("K89").Value = 6754.375
$ws.Range("L89").Value = 23717.5
$ws.Range("M89").Value = -1138.375
$ws.Range("N89").Value = -34949.5

$ws.Range("H98").Value = 1094.2
$ws.Range("I98").Value = 801.0769
$ws.Range("K98").Value = 801.0769
$ws.Range("M98").Value = 696.9231

$ws.Range("H122").Value = 1094.2
$ws.Range("I122").Value = 801.0769
$ws.Range("K122").Value = 2403.2307
$ws.Range("M122").Value = 46.76929999999993

$ws.Range("H138").Value = 6018.4585
$ws.Range("J138").Value = 6939.7334
$ws.Range("L138").Value = 20819.2002
$ws.Range("N138").Value = -31099.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 13749.25
$ws.Range("I10").Value = 1997
$ws.Range("J10").Value = 17666.666
$ws.Range("K10").Value = 1997
$ws.Range("L10").Value = 17666.666
$ws.Range("M10").Value = -1827
$ws.Range("N10").Value = -18006.666

$ws.Range("H13").Value = 3080
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 3080
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 3080
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -3368

$ws.Range("H29").Value = 2200
$ws.Range("J29").Value = 2400
$ws.Range("L29").Value = 2400
$ws.Range("N29").Value = -3016

$ws.Range("H30").Value = 1499.5
$ws.Range("I30").Value = 1499
$ws.Range("K30").Value = 1499
$ws.Range("M30").Value = -1349

$ws.Range("H122").Value = 1849.5
$ws.Range("I122").Value = 1199
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3597
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1147
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2714772.2
$ws.Range("I7").Value = 3167217
$ws.Range("J7").Value = 104
$ws.Range("K7").Value = 3167217
$ws.Range("L7").Value = 104
$ws.Range("M7").Value = -3167104
$ws.Range("N7").Value = -330

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 10867.533
$ws.Range("I12").Value = 786.3
$ws.Range("J12").Value = 31030
$ws.Range("K12").Value = 786.3
$ws.Range("L12").Value = 31030
$ws.Range("M12").Value = -616.3
$ws.Range("N12").Value = -31370

$ws.Range("H19").Value = 15462.857
$ws.Range("I19").Value = 745
$ws.Range("J19").Value = 21350
$ws.Range("K19").Value = 745
$ws.Range("L19").Value = 21350
$ws.Range("M19").Value = -575
$ws.Range("N19").Value = -21690

$ws.Range("H23").Value = 45738.145
$ws.Range("I23").Value = 42503.5
$ws.Range("J23").Value = 47032
$ws.Range("K23").Value = 42503.5
$ws.Range("L23").Value = 47032
$ws.Range("M23").Value = -42263.5
$ws.Range("N23").Value = -47512

$ws.Range("H24").Value = 15462.857
$ws.Range("I24").Value = 745
$ws.Range("J24").Value = 21350
$ws.Range("K24").Value = 745
$ws.Range("L24").Value = 21350
$ws.Range("M24").Value = -575
$ws.Range("N24").Value = -21690

$ws.Range("H27").Value = 45738.145
$ws.Range("I27").Value = 42503.5
$ws.Range("J27").Value = 47032
$ws.Range("K27").Value = 42503.5
$ws.Range("L27").Value = 47032
$ws.Range("M27").Value = -42311.5
$ws.Range("N27").Value = -47416

$ws.Range("H58").Value = 1100
$ws.Range("I58").Value = 966.6667
$ws.Range("K58").Value = 966.6667
$ws.Range("M58").Value = -763.6667

$ws.Range("H136").Value = 1100
$ws.Range("I136").Value = 966.6667
$ws.Range("K136").Value = 2900.0001
$ws.Range("M136").Value = -350.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 456
$ws.Range("J32").Value = 456
$ws.Range("L32").Value = 1368
$ws.Range("N32").Value = -1934

$ws.Range("H131").Value = 984
$ws.Range("J131").Value = 986.89655
$ws.Range("L131").Value = 2960.68965
$ws.Range("N131").Value = -13040.68965

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 9762.143
$ws.Range("I3").Value = 9444.777
$ws.Range("K3").Value = 9444.777
$ws.Range("M3").Value = -9328.777

$ws.Range("H10").Value = 253137.5
$ws.Range("J10").Value = 250
$ws.Range("L10").Value = 250
$ws.Range("N10").Value = -588

$ws.Range("H11").Value = 6500500
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 6500500
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 6500500
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -6500778

$ws.Range("H14").Value = 167785.72
$ws.Range("I14").Value = 1000000
$ws.Range("K14").Value = 1000000
$ws.Range("M14").Value = -999832

$ws.Range("H132").Value = 3847.5454
$ws.Range("I132").Value = 3462.4
$ws.Range("K132").Value = 10387.2
$ws.Range("M132").Value = -7857.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 50000
$ws.Range("J2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("N2").Value = -50224

$ws.Range("H5").Value = 53877.25
$ws.Range("I5").Value = 58503
$ws.Range("K5").Value = 58503
$ws.Range("M5").Value = -58390

$ws.Range("H12").Value = 17250
$ws.Range("J12").Value = 17250
$ws.Range("L12").Value = 17250
$ws.Range("N12").Value = -17590

$ws.Range("H22").Value = 4922.615
$ws.Range("I22").Value = 3798
$ws.Range("J22").Value = 5422.4443
$ws.Range("K22").Value = 3798
$ws.Range("L22").Value = 5422.4443
$ws.Range("M22").Value = -3503
$ws.Range("N22").Value = -6012.4443

$ws.Range("H27").Value = 4922.615
$ws.Range("I27").Value = 3798
$ws.Range("J27").Value = 5422.4443
$ws.Range("K27").Value = 3798
$ws.Range("L27").Value = 5422.4443
$ws.Range("M27").Value = -3691
$ws.Range("N27").Value = -5636.4443

$ws.Range("H68").Value = 4250
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4250
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4250
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -5748

$ws.Range("H71").Value = 4250
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4250
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 21250
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -28738

$ws.Range("H122").Value = 1994
$ws.Range("I122").Value = 1994
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5982
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3532
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4151.5
$ws.Range("I132").Value = 2960.2
$ws.Range("J132").Value = 5342.8
$ws.Range("K132").Value = 8880.599999999999
$ws.Range("L132").Value = 16028.4
$ws.Range("M132").Value = -6350.599999999999
$ws.Range("N132").Value = -21088.4

$ws.Range("H136").Value = 7953
$ws.Range("I136").Value = 7948.3
$ws.Range("K136").Value = 23844.9
$ws.Range("M136").Value = -21294.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4239.6
$ws.Range("J4").Value = 2799.5
$ws.Range("L4").Value = 2799.5
$ws.Range("N4").Value = -3025.5

$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H132").Value = 2284.3333
$ws.Range("I132").Value = 1990.75
$ws.Range("K132").Value = 5972.25
$ws.Range("M132").Value = -3442.25

$ws.Range("H136").Value = 1165.3125
$ws.Range("J136").Value = 1133
$ws.Range("L136").Value = 3399
$ws.Range("N136").Value = -8499
